# 07.03.2020 MC Sales Details
#
# Adds a "sales group" column (E) alongside the existing IMEI table and
# colour-codes the C (v48) / D (L23i) IMEI columns to show which units
# belong to which sales group / batch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colour constants (decimal RGB, R + G*256 + B*65536) -----------------------
$yellow       = 65535      # FFFF00
$orange       = 49407      # FFC000
$green        = 5287936    # 00B050
$blueGray60   = 14922894   # 8EB4E3 (Text2/Dk2, Lighter 40%)
$blueGray80   = 15849926   # C6D9F1 (Text2/Dk2, Lighter 60%)
$blue40       = 14136213   # 95B3D7 (Accent1, Lighter 40%)
$blue80       = 15918812   # DCE6F2 (Accent1, Lighter 60%)

# ---------------------------------------------------------------------
# New column E: sales-group labels, merged in pairs next to the D column
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "Desh"
$ws.Range("E5").Value = "Rose"
$ws.Range("E7").Value = "Jilani"
$ws.Range("E8").Value = "Natore"

[void]$ws.Range("E3:E4").Merge()
[void]$ws.Range("E5:E6").Merge()
[void]$ws.Range("E8:E9").Merge()

# A cell that already carries the thin boxed border used throughout the
# C/D columns -- copying its format onto the new cells lets them share
# the exact same border definition instead of creating near-duplicates.
[void]$ws.Range("C13").Copy()
[void]$ws.Range("E3:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# E column cells share the boxed-cell look used throughout columns C/D
$ws.Range("E3:E4").Interior.Color = $blueGray80
$ws.Range("E5:E6").Interior.Color = $blueGray60
$ws.Range("E7").Interior.Color = $blue80
$ws.Range("E8:E9").Interior.Color = $blue40

# ---------------------------------------------------------------------
# Column C highlight groups (v48 IMEI numbers)
# ---------------------------------------------------------------------
$ws.Range("C3:C12").Interior.Color = $yellow
$ws.Range("C13:C22").Interior.Color = $green
$ws.Range("C23:C32").Interior.Color = $orange

# ---------------------------------------------------------------------
# Column D highlight groups (L23i IMEI numbers), paired with column E
# ---------------------------------------------------------------------
$ws.Range("D3:D4").Interior.Color = $blueGray80
$ws.Range("D5:D6").Interior.Color = $blueGray60
$ws.Range("D7").Interior.Color = $blue80
$ws.Range("D8:D9").Interior.Color = $blue40

# D2 (group header, green) keeps its box but loses the bottom edge so it
# visually joins with D3 below it.
$ws.Range("D2").Borders.LineStyle = -4142
$ws.Range("D2").Borders.Item(7).LineStyle = 1
$ws.Range("D2").Borders.Item(7).Weight = 2
$ws.Range("D2").Borders.Item(10).LineStyle = 1
$ws.Range("D2").Borders.Item(10).Weight = 2
$ws.Range("D2").Borders.Item(8).LineStyle = 1
$ws.Range("D2").Borders.Item(8).Weight = 2
$ws.Range("D2").Interior.Color = $green

# D10 (closing row) keeps the left/right box edges plus a bottom edge,
# with no top edge and no fill, so the group box visually closes here.
$ws.Range("D10").Borders.LineStyle = -4142
$ws.Range("D10").Borders.Item(7).LineStyle = 1
$ws.Range("D10").Borders.Item(7).Weight = 2
$ws.Range("D10").Borders.Item(10).LineStyle = 1
$ws.Range("D10").Borders.Item(10).Weight = 2
$ws.Range("D10").Borders.Item(9).LineStyle = 1
$ws.Range("D10").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------
# Sheet view / dimension bookkeeping
# ---------------------------------------------------------------------
[void]$ws.Range("M11").Select()

Write-Host "done"
